# Update countries & provincias Spain
# Refresh the COVID data snapshot: update case counters for a handful of
# countries, update the "last updated" timestamp, and re-rank the country
# rows whose totals changed enough to swap relative order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain numeric refreshes (no change of rank) ---------------------------

# Row 4: Estados Unidos
$ws.Range("B4").Value = 1794153
$ws.Range("C4").Value = 623
$ws.Range("E4").Value = 1169992
$ws.Range("G4").Value = 8
$ws.Range("H4").Value = 104550

# Row 24: Paises Bajos
$ws.Range("B24").Value = 46257
$ws.Range("C24").Value = 131
$ws.Range("G24").Value = 20
$ws.Range("H24").Value = 5951

# Row 26: Bielorrusia
$ws.Range("B26").Value = 41658
$ws.Range("C26").Value = 894
$ws.Range("D26").Value = 17964
$ws.Range("E26").Value = 23465
$ws.Range("G26").Value = 5
$ws.Range("H26").Value = 229

# Row 28: Suecia
$ws.Range("B28").Value = 37113
$ws.Range("C28").Value = 637
$ws.Range("E28").Value = 27747
$ws.Range("G28").Value = 45
$ws.Range("H28").Value = 4395

# Row 50: Dinamarca
$ws.Range("B50").Value = 11633
$ws.Range("C50").Value = 40
$ws.Range("D50").Value = 10327
$ws.Range("E50").Value = 735
$ws.Range("G50").Value = 3
$ws.Range("H50").Value = 571

# Row 89: Croacia
$ws.Range("B89").Value = 2246
$ws.Range("C89").Value = 1
$ws.Range("D89").Value = 2063
$ws.Range("E89").Value = 80

# Row 90: Republica de Macedonia
$ws.Range("B90").Value = 2164
$ws.Range("C90").Value = 35
$ws.Range("D90").Value = 1535
$ws.Range("E90").Value = 498
$ws.Range("G90").Value = 5
$ws.Range("H90").Value = 131

# --- Re-ranked pairs (country whose total overtakes its neighbour swaps
#     places with it) ---------------------------------------------------

# Gibraltar overtakes Guadalupe (rows 162/163)
$ws.Range("A162").Value = "Gibraltar"
$ws.Range("B162").Value = 169
$ws.Range("C162").Value = 8
$ws.Range("D162").Value = 149
$ws.Range("E162").Value = 20
$ws.Range("F162").Value = 0
$ws.Range("G162").Value = 0
$ws.Range("H162").Value = 0

$ws.Range("A163").Value = "Guadalupe"
$ws.Range("B163").Value = 162
$ws.Range("C163").Value = 0
$ws.Range("D163").Value = 138
$ws.Range("E163").Value = 10
$ws.Range("F163").Value = 0
$ws.Range("G163").Value = 0
$ws.Range("H163").Value = 14

# Curazao overtakes Fiyi (rows 198/199)
$ws.Range("A198").Value = "Curazao"
$ws.Range("B198").Value = 18
$ws.Range("C198").Value = 0
$ws.Range("D198").Value = 14
$ws.Range("E198").Value = 3
$ws.Range("F198").Value = 0
$ws.Range("G198").Value = 0
$ws.Range("H198").Value = 1

$ws.Range("A199").Value = "Fiyi"
$ws.Range("B199").Value = 18
$ws.Range("C199").Value = 0
$ws.Range("D199").Value = 15
$ws.Range("E199").Value = 3
$ws.Range("F199").Value = 0
$ws.Range("G199").Value = 0
$ws.Range("H199").Value = 0

# Montserrat overtakes Seychelles (rows 210/211)
$ws.Range("A210").Value = "Montserrat"
$ws.Range("B210").Value = 11
$ws.Range("C210").Value = 0
$ws.Range("D210").Value = 10
$ws.Range("E210").Value = 0
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 0
$ws.Range("H210").Value = 1

$ws.Range("A211").Value = "Seychelles"
$ws.Range("B211").Value = 11
$ws.Range("C211").Value = 0
$ws.Range("D211").Value = 11
$ws.Range("E211").Value = 0
$ws.Range("F211").Value = 0
$ws.Range("G211").Value = 0
$ws.Range("H211").Value = 0

# Islas Virgenes Britanicas overtakes Papua Nueva Guinea (rows 213/214)
$ws.Range("A213").Value = "Islas Virgenes Britanicas"
$ws.Range("B213").Value = 8
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 7
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 1

$ws.Range("A214").Value = "Papua Nueva Guinea"
$ws.Range("B214").Value = 8
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 8
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 0

# --- Timestamp footer (title cell, row 1) -------------------------------

$ws.Range("A1").Value = "Datos actualizados a 30 de Mayo de 2020 a las 14:05"
